$wb = $excel.ActiveWorkbook

# --- OFF sheet: row 3 (R) updates ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 174
$wsOff.Range("C3").Value = 128
$wsOff.Range("D3").Value = 39
$wsOff.Range("E3").Value = 29
$wsOff.Range("F3").Value = 3
$wsOff.Range("G3").Value = 3

# --- DEF sheet: row 3 (R) updates ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 331
$wsDef.Range("C3").Value = 245
$wsDef.Range("D3").Value = 81
$wsDef.Range("E3").Value = 37
$wsDef.Range("F3").Value = 4
$wsDef.Range("G3").Value = 2
